$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.511.50'
$ws.Range("E2").Value = '  +0.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.912.11'
$ws.Range("E3").Value = '  +0.33%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.66%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.32'
$ws.Range("E5").Value = '  +0.68%  '

$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4848'
$ws.Range("E7").Value = '  +2.62%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4066'
$ws.Range("E8").Value = '  +0.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08166'
$ws.Range("E9").Value = '  +1.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.014'
$ws.Range("E10").Value = '  +2.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.49'
$ws.Range("E11").Value = '  +3.70%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.911.83'
$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.036'
$ws.Range("E13").Value = '  +3.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.201'
$ws.Range("E14").Value = '  +2.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.55'
$ws.Range("E15").Value = '  +1.37%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.008'
$ws.Range("E16").Value = '  +0.63%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06775'
$ws.Range("E17").Value = '  +2.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001038'
$ws.Range("E18").Value = '  +1.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.68'
$ws.Range("E19").Value = '  +0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.007'
$ws.Range("E20").Value = '  +0.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.527.43'
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.629'
$ws.Range("E22").Value = '  +2.33%  '

$ws.Range("E23").Value = '  +2.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.193'
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.189.65'
$ws.Range("E25").Value = '  +2.84%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.572'
$ws.Range("E26").Value = '  +8.99%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.89'
$ws.Range("E27").Value = '  +1.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.10'
$ws.Range("E28").Value = '  +1.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.117'
$ws.Range("E29").Value = '  +1.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.51'
$ws.Range("E30").Value = '  +2.26%  '

$ws.Range("E31").Value = '  -3.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09515'
$ws.Range("E32").Value = '  +0.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.527'
$ws.Range("E33").Value = '  +3.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.561'
$ws.Range("E34").Value = '  +0.46%  '

$ws.Range("E35").Value = '  -1.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02277'
$ws.Range("E36").Value = '  +1.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06119'
$ws.Range("E37").Value = '  +0.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.180'
$ws.Range("E38").Value = '  +0.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.86'
$ws.Range("E39").Value = '  +8.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5966'
$ws.Range("E40").Value = '  +2.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.969'
$ws.Range("E41").Value = '  -1.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1859'
$ws.Range("E42").Value = '  +1.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.281'
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.376'
$ws.Range("E44").Value = '  -4.82%  '

$ws.Range("E45").Value = '  +2.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07624'
$ws.Range("E46").Value = '  -2.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5584'
$ws.Range("E47").Value = '  +1.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.949'
$ws.Range("E48").Value = '  +2.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.66'
$ws.Range("E49").Value = '  +3.10%  '

$ws.Range("E50").Value = '  +2.29%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.416'
$ws.Range("E51").Value = '  +3.10%  '
